# Remove rows referencing the custom "annotated" SARS1/SARS2 point-mutant
# constructs (SARS2+Q498Y, SARS2+N501Y, SARS2+Q498Y+N501Y, SARS1+Y498Q,
# SARS1+T501Y, SARS1+Y498Q+T501Y) and close up a blank spacer row higher
# in the sheet, per the commit "remove custom annotated pdfs and ppts".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the six rows holding the custom annotated variant sequences.
$ws.Rows("86:91").Delete()

# Close up one of the blank rows between the HuB2013 row (74) and the
# Rc-o319 row (formerly 77, now 76).
$ws.Rows("76:76").Delete()

# Reflect the user's final selection/scroll position in the sheet.
$ws.Range("D14").Select()
